# Update column F ("dSF") values for specific rows per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    2  = 0
    3  = -2
    4  = -4
    7  = -6
    8  = -1
    9  = -6
    13 = -3
    16 = 4
    17 = 8
    18 = -4
    20 = 0
    23 = -1
    24 = -3
    28 = 2
    31 = -6
    32 = -3
    33 = 0
    34 = 1
    40 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
